# Auto-generated Excel COM-interop script to apply the Chocobo_Profits price/profit update.
# For each (sheet,row) pair identified from the diff, verify the "Leve Item ID" (column G)
# still matches the expected value before writing new currentAveragePrice / LevePrice / LeveProfit figures.

$wb = $excel.ActiveWorkbook
$mismatches = 0

$alc = $wb.Worksheets.Item("ALC")
# ALC!row 2 (Leve Item ID = 5489)
if ($alc.Range("G2").Value2 -ne 5489) {
    Write-Host "MISMATCH at ALC!G2: expected 5489, found $($alc.Range('G2').Value2)"
    $mismatches++
} else {
    $alc.Range("H2").Value = 563.5
    $alc.Range("I2").Value = 1000.5
    $alc.Range("J2").Value = 345
    $alc.Range("K2").Value = 1000.5
    $alc.Range("L2").Value = 345
    $alc.Range("M2").Value = -887.5
    $alc.Range("N2").Value = -571
}

# ALC!row 62 (Leve Item ID = 27781)
if ($alc.Range("G62").Value2 -ne 27781) {
    Write-Host "MISMATCH at ALC!G62: expected 27781, found $($alc.Range('G62').Value2)"
    $mismatches++
} else {
    $alc.Range("H62").Value = 4023.4285
    $alc.Range("I62").Value = 2194
    $alc.Range("K62").Value = 2194
    $alc.Range("M62").Value = -1570
}

# ALC!row 65 (Leve Item ID = 27781)
if ($alc.Range("G65").Value2 -ne 27781) {
    Write-Host "MISMATCH at ALC!G65: expected 27781, found $($alc.Range('G65').Value2)"
    $mismatches++
} else {
    $alc.Range("H65").Value = 4023.4285
    $alc.Range("I65").Value = 2194
    $alc.Range("K65").Value = 10970
    $alc.Range("M65").Value = -7850
}

# ALC!row 107 (Leve Item ID = 27766)
if ($alc.Range("G107").Value2 -ne 27766) {
    Write-Host "MISMATCH at ALC!G107: expected 27766, found $($alc.Range('G107').Value2)"
    $mismatches++
} else {
    $alc.Range("H107").Value = 1875.1177
    $alc.Range("I107").Value = 2641.875
    $alc.Range("J107").Value = 1193.5555
    $alc.Range("K107").Value = 2641.875
    $alc.Range("L107").Value = 1193.5555
    $alc.Range("M107").Value = -721.875
    $alc.Range("N107").Value = -5033.5555
}

# ALC!row 127 (Leve Item ID = 36114)
if ($alc.Range("G127").Value2 -ne 36114) {
    Write-Host "MISMATCH at ALC!G127: expected 36114, found $($alc.Range('G127').Value2)"
    $mismatches++
} else {
    $alc.Range("H127").Value = 720.44446
    $alc.Range("I127").Value = 460.5
    $alc.Range("J127").Value = 2800
    $alc.Range("K127").Value = 1381.5
    $alc.Range("L127").Value = 8400
    $alc.Range("M127").Value = 3578.5
    $alc.Range("N127").Value = -18320
}

# ALC!row 141 (Leve Item ID = 44161)
if ($alc.Range("G141").Value2 -ne 44161) {
    Write-Host "MISMATCH at ALC!G141: expected 44161, found $($alc.Range('G141').Value2)"
    $mismatches++
} else {
    $alc.Range("H141").Value = 65139.625
    $alc.Range("I141").Value = 79179.92
    $alc.Range("K141").Value = 237539.76
    $alc.Range("M141").Value = -232359.76
}

$arm = $wb.Worksheets.Item("ARM")
# ARM!row 5 (Leve Item ID = 5091)
if ($arm.Range("G5").Value2 -ne 5091) {
    Write-Host "MISMATCH at ARM!G5: expected 5091, found $($arm.Range('G5').Value2)"
    $mismatches++
} else {
    $arm.Range("H5").Value = 381.33334
    $arm.Range("I5").Value = 381.33334
    $arm.Range("J5").Value = 0
    $arm.Range("K5").Value = 381.33334
    $arm.Range("L5").Value = 0
    $arm.Range("M5").Value = -269.33334
    $arm.Range("N5").ClearContents() | Out-Null
}

# ARM!row 32 (Leve Item ID = 44147)
if ($arm.Range("G32").Value2 -ne 44147) {
    Write-Host "MISMATCH at ARM!G32: expected 44147, found $($arm.Range('G32').Value2)"
    $mismatches++
} else {
    $arm.Range("H32").Value = 5401.017
    $arm.Range("I32").Value = 3634.4614
    $arm.Range("J32").Value = 8845.799999999999
    $arm.Range("K32").Value = 3634.4614
    $arm.Range("L32").Value = 8845.799999999999
    $arm.Range("M32").Value = -3347.4614
    $arm.Range("N32").Value = -9419.799999999999
}

$bsm = $wb.Worksheets.Item("BSM")
# BSM!row 4 (Leve Item ID = 5091)
if ($bsm.Range("G4").Value2 -ne 5091) {
    Write-Host "MISMATCH at BSM!G4: expected 5091, found $($bsm.Range('G4').Value2)"
    $mismatches++
} else {
    $bsm.Range("H4").Value = 381.33334
    $bsm.Range("I4").Value = 381.33334
    $bsm.Range("J4").Value = 0
    $bsm.Range("K4").Value = 381.33334
    $bsm.Range("L4").Value = 0
    $bsm.Range("M4").Value = -266.33334
    $bsm.Range("N4").ClearContents() | Out-Null
}

# BSM!row 94 (Leve Item ID = 19939)
if ($bsm.Range("G94").Value2 -ne 19939) {
    Write-Host "MISMATCH at BSM!G94: expected 19939, found $($bsm.Range('G94').Value2)"
    $mismatches++
} else {
    $bsm.Range("H94").Value = 1720
    $bsm.Range("I94").Value = 1720
    $bsm.Range("J94").Value = 0
    $bsm.Range("K94").Value = 1720
    $bsm.Range("L94").Value = 0
    $bsm.Range("M94").Value = -1269
    $bsm.Range("N94").ClearContents() | Out-Null
}

$crp = $wb.Worksheets.Item("CRP")
# CRP!row 7 (Leve Item ID = 5361)
if ($crp.Range("G7").Value2 -ne 5361) {
    Write-Host "MISMATCH at CRP!G7: expected 5361, found $($crp.Range('G7').Value2)"
    $mismatches++
} else {
    $crp.Range("H7").Value = 331.46667
    $crp.Range("I7").Value = 398.875
    $crp.Range("J7").Value = 254.42857
    $crp.Range("K7").Value = 398.875
    $crp.Range("L7").Value = 254.42857
    $crp.Range("M7").Value = -285.875
    $crp.Range("N7").Value = -480.42857
}

# CRP!row 16 (Leve Item ID = 27691)
if ($crp.Range("G16").Value2 -ne 27691) {
    Write-Host "MISMATCH at CRP!G16: expected 27691, found $($crp.Range('G16').Value2)"
    $mismatches++
} else {
    $crp.Range("H16").Value = 15873687
    $crp.Range("I16").Value = 22222882
    $crp.Range("J16").Value = 700
    $crp.Range("K16").Value = 22222882
    $crp.Range("L16").Value = 700
    $crp.Range("M16").Value = -22222595
    $crp.Range("N16").Value = -1274
}

# CRP!row 22 (Leve Item ID = 5367)
if ($crp.Range("G22").Value2 -ne 5367) {
    Write-Host "MISMATCH at CRP!G22: expected 5367, found $($crp.Range('G22').Value2)"
    $mismatches++
} else {
    $crp.Range("H22").Value = 847.7727
    $crp.Range("I22").Value = 496.35715
    $crp.Range("J22").Value = 1462.75
    $crp.Range("K22").Value = 496.35715
    $crp.Range("L22").Value = 1462.75
    $crp.Range("M22").Value = -146.35715
    $crp.Range("N22").Value = -2162.75
}

# CRP!row 31 (Leve Item ID = 44023)
if ($crp.Range("G31").Value2 -ne 44023) {
    Write-Host "MISMATCH at CRP!G31: expected 44023, found $($crp.Range('G31').Value2)"
    $mismatches++
} else {
    $crp.Range("H31").Value = 2496.5881
    $crp.Range("I31").Value = 1205.7826
    $crp.Range("J31").Value = 5195.5454
    $crp.Range("K31").Value = 1205.7826
    $crp.Range("L31").Value = 5195.5454
    $crp.Range("M31").Value = -910.7826
    $crp.Range("N31").Value = -5785.5454
}

# CRP!row 34 (Leve Item ID = 44023)
if ($crp.Range("G34").Value2 -ne 44023) {
    Write-Host "MISMATCH at CRP!G34: expected 44023, found $($crp.Range('G34').Value2)"
    $mismatches++
} else {
    $crp.Range("H34").Value = 2496.5881
    $crp.Range("I34").Value = 1205.7826
    $crp.Range("J34").Value = 5195.5454
    $crp.Range("K34").Value = 1205.7826
    $crp.Range("L34").Value = 5195.5454
    $crp.Range("M34").Value = -1003.7826
    $crp.Range("N34").Value = -5599.5454
}

# CRP!row 58 (Leve Item ID = 44021)
if ($crp.Range("G58").Value2 -ne 44021) {
    Write-Host "MISMATCH at CRP!G58: expected 44021, found $($crp.Range('G58').Value2)"
    $mismatches++
} else {
    $crp.Range("H58").Value = 2974.1904
    $crp.Range("I58").Value = 1805.48
    $crp.Range("J58").Value = 7469.231
    $crp.Range("K58").Value = 1805.48
    $crp.Range("L58").Value = 7469.231
    $crp.Range("M58").Value = -1602.48
    $crp.Range("N58").Value = -7875.231
}

# CRP!row 113 (Leve Item ID = 27691)
if ($crp.Range("G113").Value2 -ne 27691) {
    Write-Host "MISMATCH at CRP!G113: expected 27691, found $($crp.Range('G113').Value2)"
    $mismatches++
} else {
    $crp.Range("H113").Value = 15873687
    $crp.Range("I113").Value = 22222882
    $crp.Range("J113").Value = 700
    $crp.Range("K113").Value = 22222882
    $crp.Range("L113").Value = 700
    $crp.Range("M113").Value = -22220712
    $crp.Range("N113").Value = -5040
}

# CRP!row 122 (Leve Item ID = 36196)
if ($crp.Range("G122").Value2 -ne 36196) {
    Write-Host "MISMATCH at CRP!G122: expected 36196, found $($crp.Range('G122').Value2)"
    $mismatches++
} else {
    $crp.Range("H122").Value = 3810
    $crp.Range("I122").Value = 1515
    $crp.Range("K122").Value = 4545
    $crp.Range("M122").Value = -2095
}

# CRP!row 132 (Leve Item ID = 44019)
if ($crp.Range("G132").Value2 -ne 44019) {
    Write-Host "MISMATCH at CRP!G132: expected 44019, found $($crp.Range('G132').Value2)"
    $mismatches++
} else {
    $crp.Range("H132").Value = 2660.5122
    $crp.Range("I132").Value = 2184.8286
    $crp.Range("K132").Value = 6554.485799999999
    $crp.Range("M132").Value = -4024.485799999999
}

# CRP!row 134 (Leve Item ID = 44020)
if ($crp.Range("G134").Value2 -ne 44020) {
    Write-Host "MISMATCH at CRP!G134: expected 44020, found $($crp.Range('G134').Value2)"
    $mismatches++
} else {
    $crp.Range("H134").Value = 2047.3214
    $crp.Range("I134").Value = 1330.55
    $crp.Range("J134").Value = 3839.25
    $crp.Range("K134").Value = 3991.65
    $crp.Range("L134").Value = 11517.75
    $crp.Range("M134").Value = -1456.65
    $crp.Range("N134").Value = -16587.75
}

# CRP!row 135 (Leve Item ID = 42008)
if ($crp.Range("G135").Value2 -ne 42008) {
    Write-Host "MISMATCH at CRP!G135: expected 42008, found $($crp.Range('G135').Value2)"
    $mismatches++
} else {
    $crp.Range("H135").Value = 39832
    $crp.Range("J135").Value = 39832
    $crp.Range("L135").Value = 39832
    $crp.Range("N135").Value = -49972
}

# CRP!row 136 (Leve Item ID = 44021)
if ($crp.Range("G136").Value2 -ne 44021) {
    Write-Host "MISMATCH at CRP!G136: expected 44021, found $($crp.Range('G136').Value2)"
    $mismatches++
} else {
    $crp.Range("H136").Value = 2974.1904
    $crp.Range("I136").Value = 1805.48
    $crp.Range("J136").Value = 7469.231
    $crp.Range("K136").Value = 5416.440000000001
    $crp.Range("L136").Value = 22407.693
    $crp.Range("M136").Value = -2866.440000000001
    $crp.Range("N136").Value = -27507.693
}

$cul = $wb.Worksheets.Item("CUL")
# CUL!row 107 (Leve Item ID = 27838)
if ($cul.Range("G107").Value2 -ne 27838) {
    Write-Host "MISMATCH at CUL!G107: expected 27838, found $($cul.Range('G107').Value2)"
    $mismatches++
} else {
    $cul.Range("H107").Value = 143655.86
    $cul.Range("J107").Value = 250800.25
    $cul.Range("L107").Value = 752400.75
    $cul.Range("N107").Value = -756240.75
}

# CUL!row 137 (Leve Item ID = 44088)
if ($cul.Range("G137").Value2 -ne 44088) {
    Write-Host "MISMATCH at CUL!G137: expected 44088, found $($cul.Range('G137').Value2)"
    $mismatches++
} else {
    $cul.Range("H137").Value = 945
    $cul.Range("I137").Value = 945
    $cul.Range("K137").Value = 2835
    $cul.Range("M137").Value = 2265
}

# CUL!row 140 (Leve Item ID = 44097)
if ($cul.Range("G140").Value2 -ne 44097) {
    Write-Host "MISMATCH at CUL!G140: expected 44097, found $($cul.Range('G140').Value2)"
    $mismatches++
} else {
    $cul.Range("H140").Value = 2753.7368
    $cul.Range("I140").Value = 2905.1765
    $cul.Range("J140").Value = 1466.5
    $cul.Range("K140").Value = 8715.529500000001
    $cul.Range("L140").Value = 4399.5
    $cul.Range("M140").Value = -3535.529500000001
    $cul.Range("N140").Value = -14759.5
}

$gsm = $wb.Worksheets.Item("GSM")
# GSM!row 2 (Leve Item ID = 5062)
if ($gsm.Range("G2").Value2 -ne 5062) {
    Write-Host "MISMATCH at GSM!G2: expected 5062, found $($gsm.Range('G2').Value2)"
    $mismatches++
} else {
    $gsm.Range("H2").Value = 211.92857
    $gsm.Range("I2").Value = 63
    $gsm.Range("J2").Value = 480
    $gsm.Range("K2").Value = 63
    $gsm.Range("L2").Value = 480
    $gsm.Range("M2").Value = 50
    $gsm.Range("N2").Value = -706
}

# GSM!row 11 (Leve Item ID = 4422)
if ($gsm.Range("G11").Value2 -ne 4422) {
    Write-Host "MISMATCH at GSM!G11: expected 4422, found $($gsm.Range('G11').Value2)"
    $mismatches++
} else {
    $gsm.Range("H11").Value = 7217834.5
    $gsm.Range("I11").Value = 12857571
    $gsm.Range("K11").Value = 12857571
    $gsm.Range("M11").Value = -12857432
}

# GSM!row 70 (Leve Item ID = 14146)
if ($gsm.Range("G70").Value2 -ne 14146) {
    Write-Host "MISMATCH at GSM!G70: expected 14146, found $($gsm.Range('G70').Value2)"
    $mismatches++
} else {
    $gsm.Range("H70").Value = 5656.9614
    $gsm.Range("I70").Value = 5262.5835
    $gsm.Range("J70").Value = 6544.3125
    $gsm.Range("K70").Value = 5262.5835
    $gsm.Range("L70").Value = 6544.3125
    $gsm.Range("M70").Value = -4992.5835
    $gsm.Range("N70").Value = -7084.3125
}

# GSM!row 73 (Leve Item ID = 14146)
if ($gsm.Range("G73").Value2 -ne 14146) {
    Write-Host "MISMATCH at GSM!G73: expected 14146, found $($gsm.Range('G73').Value2)"
    $mismatches++
} else {
    $gsm.Range("H73").Value = 5656.9614
    $gsm.Range("I73").Value = 5262.5835
    $gsm.Range("J73").Value = 6544.3125
    $gsm.Range("K73").Value = 5262.5835
    $gsm.Range("L73").Value = 6544.3125
    $gsm.Range("M73").Value = -4326.5835
    $gsm.Range("N73").Value = -8416.3125
}

# GSM!row 124 (Leve Item ID = 34247)
if ($gsm.Range("G124").Value2 -ne 34247) {
    Write-Host "MISMATCH at GSM!G124: expected 34247, found $($gsm.Range('G124').Value2)"
    $mismatches++
} else {
    $gsm.Range("H124").Value = 41572.715
    $gsm.Range("J124").Value = 41572.715
    $gsm.Range("L124").Value = 41572.715
    $gsm.Range("N124").Value = -51392.715
}

# GSM!row 126 (Leve Item ID = 36184)
if ($gsm.Range("G126").Value2 -ne 36184) {
    Write-Host "MISMATCH at GSM!G126: expected 36184, found $($gsm.Range('G126').Value2)"
    $mismatches++
} else {
    $gsm.Range("H126").Value = 3336.6262
    $gsm.Range("I126").Value = 2931.0676
    $gsm.Range("J126").Value = 4537.08
    $gsm.Range("K126").Value = 8793.202799999999
    $gsm.Range("L126").Value = 13611.24
    $gsm.Range("M126").Value = -6323.202799999999
    $gsm.Range("N126").Value = -18551.24
}

# GSM!row 132 (Leve Item ID = 44008)
if ($gsm.Range("G132").Value2 -ne 44008) {
    Write-Host "MISMATCH at GSM!G132: expected 44008, found $($gsm.Range('G132').Value2)"
    $mismatches++
} else {
    $gsm.Range("H132").Value = 3700
    $gsm.Range("I132").Value = 2646.8
    $gsm.Range("J132").Value = 6859.6
    $gsm.Range("K132").Value = 7940.400000000001
    $gsm.Range("L132").Value = 20578.8
    $gsm.Range("M132").Value = -5410.400000000001
    $gsm.Range("N132").Value = -25638.8
}

# GSM!row 133 (Leve Item ID = 41854)
if ($gsm.Range("G133").Value2 -ne 41854) {
    Write-Host "MISMATCH at GSM!G133: expected 41854, found $($gsm.Range('G133').Value2)"
    $mismatches++
} else {
    $gsm.Range("H133").Value = 40711.668
    $gsm.Range("J133").Value = 40711.668
    $gsm.Range("L133").Value = 40711.668
    $gsm.Range("N133").Value = -50831.668
}

$ltw = $wb.Worksheets.Item("LTW")
# LTW!row 16 (Leve Item ID = 5289)
if ($ltw.Range("G16").Value2 -ne 5289) {
    Write-Host "MISMATCH at LTW!G16: expected 5289, found $($ltw.Range('G16').Value2)"
    $mismatches++
} else {
    $ltw.Range("H16").Value = 979.5
    $ltw.Range("I16").Value = 979.5
    $ltw.Range("K16").Value = 979.5
    $ltw.Range("M16").Value = -809.5
}

# LTW!row 88 (Leve Item ID = 10961)
if ($ltw.Range("G88").Value2 -ne 10961) {
    Write-Host "MISMATCH at LTW!G88: expected 10961, found $($ltw.Range('G88').Value2)"
    $mismatches++
} else {
    $ltw.Range("H88").Value = 2000
    $ltw.Range("I88").Value = 2000
    $ltw.Range("K88").Value = 2000
    $ltw.Range("M88").Value = -1572
}

# LTW!row 91 (Leve Item ID = 10961)
if ($ltw.Range("G91").Value2 -ne 10961) {
    Write-Host "MISMATCH at LTW!G91: expected 10961, found $($ltw.Range('G91').Value2)"
    $mismatches++
} else {
    $ltw.Range("H91").Value = 2000
    $ltw.Range("I91").Value = 2000
    $ltw.Range("K91").Value = 2000
    $ltw.Range("M91").Value = -518
}

# LTW!row 100 (Leve Item ID = 19995)
if ($ltw.Range("G100").Value2 -ne 19995) {
    Write-Host "MISMATCH at LTW!G100: expected 19995, found $($ltw.Range('G100').Value2)"
    $mismatches++
} else {
    $ltw.Range("H100").Value = 2726.7368
    $ltw.Range("I100").Value = 1557.1428
    $ltw.Range("J100").Value = 3409
    $ltw.Range("K100").Value = 1557.1428
    $ltw.Range("L100").Value = 3409
    $ltw.Range("M100").Value = -1016.1428
    $ltw.Range("N100").Value = -4491
}

# LTW!row 103 (Leve Item ID = 18526)
if ($ltw.Range("G103").Value2 -ne 18526) {
    Write-Host "MISMATCH at LTW!G103: expected 18526, found $($ltw.Range('G103').Value2)"
    $mismatches++
} else {
    $ltw.Range("H103").Value = 0
    $ltw.Range("J103").Value = 0
    $ltw.Range("L103").Value = 0
    $ltw.Range("N103").ClearContents() | Out-Null
}

# LTW!row 136 (Leve Item ID = 44060)
if ($ltw.Range("G136").Value2 -ne 44060) {
    Write-Host "MISMATCH at LTW!G136: expected 44060, found $($ltw.Range('G136').Value2)"
    $mismatches++
} else {
    $ltw.Range("H136").Value = 6017.2666
    $ltw.Range("J136").Value = 9033.333000000001
    $ltw.Range("L136").Value = 27099.999
    $ltw.Range("N136").Value = -32199.999
}

$wvr = $wb.Worksheets.Item("WVR")
# WVR!row 132 (Leve Item ID = 44029)
if ($wvr.Range("G132").Value2 -ne 44029) {
    Write-Host "MISMATCH at WVR!G132: expected 44029, found $($wvr.Range('G132').Value2)"
    $mismatches++
} else {
    $wvr.Range("H132").Value = 15158361
    $wvr.Range("I132").Value = 7181.7646
    $wvr.Range("K132").Value = 21545.2938
    $wvr.Range("M132").Value = -19015.2938
}

if ($mismatches -eq 0) {
    Write-Host "All 40 target rows updated successfully."
} else {
    Write-Host "$mismatches row(s) failed the Leve Item ID sanity check."
}
